# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-safe assignments (Coin names, Links, Volume %, multi-dot prices) ---
$ws.Range('D2').Value = '25.765.67'
$ws.Range('E2').Value = '  -3.75%  '
$ws.Range('D3').Value = '1.819.63'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -7.33%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('E7').Value = '  -4.12%  '
$ws.Range('E8').Value = '  -5.23%  '
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('E10').Value = '  -6.86%  '
$ws.Range('E11').Value = '  -7.51%  '
$ws.Range('E12').Value = '  -6.49%  '
$ws.Range('E13').Value = '  -3.40%  '
$ws.Range('D14').Value = '1.807.01'
$ws.Range('E14').Value = '  -6.64%  '
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('E16').Value = '  -5.34%  '
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -4.58%  '
$ws.Range('E19').Value = '  -5.15%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '25.807.35'
$ws.Range('E21').Value = '  -3.72%  '
$ws.Range('E22').Value = '  -4.40%  '
$ws.Range('E23').Value = '  -5.72%  '
$ws.Range('E24').Value = '  -3.93%  '
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('E26').Value = '  -2.60%  '
$ws.Range('E27').Value = '  -3.73%  '
$ws.Range('E28').Value = '  -4.70%  '
$ws.Range('E29').Value = '  -3.89%  '
$ws.Range('E30').Value = '  -7.87%  '
$ws.Range('E31').Value = '  -8.03%  '
$ws.Range('E32').Value = '  -3.92%  '
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('E34').Value = '  -9.11%  '
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('E36').Value = '  -2.19%  '
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E38').Value = '  -11.47%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E39').Value = '  -4.54%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E40').Value = '  -15.01%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E41').Value = '  -9.18%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E43').Value = '  -4.05%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E44').Value = '  -8.20%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E46').Value = '  -11.97%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E47').Value = '  -8.01%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E49').Value = '  -6.79%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E50').Value = '  -8.40%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E51').Value = '  -3.31%  '

# --- Price values that Excel would otherwise coerce to numbers; force text ---
$numericPriceCells = @{
    'D5' = '278.64'
    'D6' = '1.001'
    'D7' = '0.5106'
    'D8' = '0.3548'
    'D9' = '44.48'
    'D10' = '0.06688'
    'D11' = '20.01'
    'D12' = '0.8304'
    'D13' = '0.07871'
    'D15' = '5.089'
    'D16' = '87.89'
    'D18' = '14.13'
    'D19' = '0.000008047'
    'D20' = '1.001'
    'D22' = '4.759'
    'D24' = '6.128'
    'D26' = '142.25'
    'D27' = '1.671'
    'D28' = '17.16'
    'D29' = '109.44'
    'D30' = '4.343'
    'D31' = '4.244'
    'D32' = '0.08774'
    'D33' = '0.04911'
    'D34' = '0.7305'
    'D35' = '1.137'
    'D37' = '3.157'
    'D38' = '2.362'
    'D39' = '0.01861'
    'D40' = '0.5190'
    'D41' = '0.9679'
    'D42' = '115.10'
    'D43' = '6.236'
    'D44' = '8.027'
    'D45' = '1.000'
    'D46' = '0.4555'
    'D47' = '0.1371'
    'D48' = '36.56'
    'D49' = '9.267'
    'D50' = '1.502'
    'D51' = '0.05837'
}
foreach ($cell in $numericPriceCells.Keys) {
    $ws.Range($cell).NumberFormat = "@"
}
foreach ($cell in $numericPriceCells.Keys) {
    $ws.Range($cell).Value = $numericPriceCells[$cell]
}
foreach ($cell in $numericPriceCells.Keys) {
    $ws.Range($cell).ClearFormats()
}
